$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 313 (pushes old 313-318 down to 315-320,
# preserving their values/format untouched).
$ws.Rows("313:314").Insert()

# Row 313 - new weekly entry
$ws.Range("A313").Value = 10
$ws.Range("B313").Value = "Vega Modelo de Temuco"
$ws.Range("C313").Value = "La Araucanía"
$ws.Range("D313").Value = 44448
$ws.Range("E313").Value = 9
$ws.Range("F313").Value = "Fruta"
$ws.Range("G313").Value = 100108
$ws.Range("H313").Value = "Tropicales y subtropicales"
$ws.Range("I313").Value = 100108006
$ws.Range("J313").Value = "Plátano"
$ws.Range("K313").Value = "Barraganete"
$ws.Range("L313").Value = "Primera"
$ws.Range("M313").Value = 40
$ws.Range("N313").Value = 24000
$ws.Range("O313").Value = 24000
$ws.Range("P313").Value = 24000
$ws.Range("Q313").Value = '$/caja 20 kilos'
$ws.Range("R313").Value = "Ecuador"
$ws.Range("S313").Value = 1200
$ws.Range("T313").Value = 20

# Row 314 - new weekly entry
$ws.Range("A314").Value = 10
$ws.Range("B314").Value = "Vega Modelo de Temuco"
$ws.Range("C314").Value = "La Araucanía"
$ws.Range("D314").Value = 44448
$ws.Range("E314").Value = 9
$ws.Range("F314").Value = "Fruta"
$ws.Range("G314").Value = 100108
$ws.Range("H314").Value = "Tropicales y subtropicales"
$ws.Range("I314").Value = 100108006
$ws.Range("J314").Value = "Plátano"
$ws.Range("K314").Value = "Sin especificar"
$ws.Range("L314").Value = "Pintón"
$ws.Range("M314").Value = 800
$ws.Range("N314").Value = 22000
$ws.Range("O314").Value = 22000
$ws.Range("P314").Value = 22000
$ws.Range("Q314").Value = '$/caja 20 kilos'
$ws.Range("R314").Value = "Ecuador"
$ws.Range("S314").Value = 1100
$ws.Range("T314").Value = 20
